$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Degree symbol (U+00B0)
$deg = [string][char]0x00B0

for ($r = 8; $r -le 15; $r++) {
    $ws.Range("G$r").Value = "CONDENSATE"
    $ws.Range("L$r").Value = "200" + $deg + "C"
    $ws.Range("M$r").Value = "1,000 kPaG"
    $ws.Range("N$r").Value = "185" + $deg + "C"
    $ws.Range("O$r").Value = "1,000 kPaG"
}

$ws.Range("E14").Value = "Pressure Retaining Bolt & Nut"
$ws.Range("E15").Value = "External Fittings"
